$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-02-23"

# Update the header label for the current (partial) month column
$ws.Range("B1").Value = "February 2022 (through February 23)"

# Update / add the daily carjacking counts for 2022-03-03's data pull
$ws.Range("D3").Value = 12     # Austin
$ws.Range("D9").Value = 4      # United Center
$ws.Range("J10").Value = 2     # Grand Crossing
$ws.Range("B11").Value = 7     # Garfield Park
$ws.Range("F13").Value = 1     # Little Italy, UIC
$ws.Range("B15").Value = 3     # West Town
$ws.Range("B20").Value = 3     # Kenwood
$ws.Range("B21").Value = 1     # Logan Square
$ws.Range("H22").Value = 1     # Humboldt Park
$ws.Range("D27").Value = 1     # Roseland
$ws.Range("D36").Value = 3     # South Chicago
$ws.Range("L46").Value = 2     # Lower West Side
$ws.Range("B54").Value = 1     # Belmont Cragin
$ws.Range("J55").Value = 1     # Avalon Park
$ws.Range("B72").Value = 2     # Lincoln Park
